# Auto-generated: apply market-data value updates scraped by the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(41, 8).Value = 294.33334  # H41
$ws.Cells.Item(41, 10).Value = 433  # J41
$ws.Cells.Item(41, 12).Value = 433  # L41
$ws.Cells.Item(41, 14).Value = -1313  # N41
$ws.Cells.Item(93, 8).Value = 601  # H93
$ws.Cells.Item(93, 10).Value = 601  # J93
$ws.Cells.Item(93, 12).Value = 601  # L93
$ws.Cells.Item(93, 14).Value = -5593  # N93
$ws.Cells.Item(116, 8).Value = 10159.4  # H116
$ws.Cells.Item(116, 9).Value = 9286.286  # I116
$ws.Cells.Item(116, 10).Value = 10923.375  # J116
$ws.Cells.Item(116, 11).Value = 9286.286  # K116
$ws.Cells.Item(116, 12).Value = 10923.375  # L116
$ws.Cells.Item(116, 13).Value = -5844.286  # M116
$ws.Cells.Item(116, 14).Value = -17807.375  # N116
$ws.Cells.Item(132, 8).Value = 10070.88  # H132
$ws.Cells.Item(132, 9).Value = 4360.45  # I132
$ws.Cells.Item(132, 10).Value = 32912.6  # J132
$ws.Cells.Item(132, 11).Value = 13081.35  # K132
$ws.Cells.Item(132, 12).Value = 98737.79999999999  # L132
$ws.Cells.Item(132, 13).Value = -10551.35  # M132
$ws.Cells.Item(132, 14).Value = -103797.8  # N132
$ws.Cells.Item(140, 8).Value = 122399  # H140
$ws.Cells.Item(140, 10).Value = 148975  # J140
$ws.Cells.Item(140, 12).Value = 148975  # L140
$ws.Cells.Item(140, 14).Value = -159335  # N140

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5816.2104  # H32
$ws.Cells.Item(32, 9).Value = 4283.7256  # I32
$ws.Cells.Item(32, 11).Value = 4283.7256  # K32
$ws.Cells.Item(32, 13).Value = -3996.7256  # M32
$ws.Cells.Item(61, 8).Value = 13999.667  # H61
$ws.Cells.Item(61, 9).Value = 2000  # I61
$ws.Cells.Item(61, 10).Value = 19999.5  # J61
$ws.Cells.Item(61, 11).Value = 2000  # K61
$ws.Cells.Item(61, 12).Value = 19999.5  # L61
$ws.Cells.Item(61, 13).Value = -1788  # M61
$ws.Cells.Item(61, 14).Value = -20423.5  # N61
$ws.Cells.Item(74, 8).Value = 23781.5  # H74
$ws.Cells.Item(74, 9).Value = 28337.334  # I74
$ws.Cells.Item(74, 10).Value = 10114  # J74
$ws.Cells.Item(74, 11).Value = 28337.334  # K74
$ws.Cells.Item(74, 12).Value = 10114  # L74
$ws.Cells.Item(74, 13).Value = -27463.334  # M74
$ws.Cells.Item(74, 14).Value = -11862  # N74
$ws.Cells.Item(77, 8).Value = 23781.5  # H77
$ws.Cells.Item(77, 9).Value = 28337.334  # I77
$ws.Cells.Item(77, 10).Value = 10114  # J77
$ws.Cells.Item(77, 11).Value = 141686.67  # K77
$ws.Cells.Item(77, 12).Value = 50570  # L77
$ws.Cells.Item(77, 13).Value = -137318.67  # M77
$ws.Cells.Item(77, 14).Value = -59306  # N77
$ws.Cells.Item(124, 8).Value = 76597.60000000001  # H124
$ws.Cells.Item(124, 10).Value = 76597.60000000001  # J124
$ws.Cells.Item(124, 12).Value = 76597.60000000001  # L124
$ws.Cells.Item(124, 14).Value = -86417.60000000001  # N124
$ws.Cells.Item(132, 8).Value = 3319.2122  # H132
$ws.Cells.Item(132, 9).Value = 3341.7188  # I132
$ws.Cells.Item(132, 11).Value = 10025.1564  # K132
$ws.Cells.Item(132, 13).Value = -7495.1564  # M132
$ws.Cells.Item(136, 8).Value = 13999.667  # H136
$ws.Cells.Item(136, 9).Value = 2000  # I136
$ws.Cells.Item(136, 10).Value = 19999.5  # J136
$ws.Cells.Item(136, 11).Value = 6000  # K136
$ws.Cells.Item(136, 12).Value = 59998.5  # L136
$ws.Cells.Item(136, 13).Value = -3450  # M136
$ws.Cells.Item(136, 14).Value = -65098.5  # N136
$ws.Cells.Item(141, 8).Value = 83810  # H141
$ws.Cells.Item(141, 10).Value = 83810  # J141
$ws.Cells.Item(141, 12).Value = 83810  # L141
$ws.Cells.Item(141, 14).Value = -94170  # N141

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 3943.2144  # H20
$ws.Cells.Item(20, 10).Value = 3876  # J20
$ws.Cells.Item(20, 12).Value = 3876  # L20
$ws.Cells.Item(20, 14).Value = -4370  # N20
$ws.Cells.Item(70, 8).Value = 169985  # H70
$ws.Cells.Item(70, 10).Value = 169985  # J70
$ws.Cells.Item(70, 12).Value = 169985  # L70
$ws.Cells.Item(70, 14).Value = -170571  # N70
$ws.Cells.Item(73, 8).Value = 169985  # H73
$ws.Cells.Item(73, 10).Value = 169985  # J73
$ws.Cells.Item(73, 12).Value = 169985  # L73
$ws.Cells.Item(73, 14).Value = -172013  # N73
$ws.Cells.Item(99, 8).Value = 9999.5  # H99
$ws.Cells.Item(99, 9).Value = 9999.5  # I99
$ws.Cells.Item(99, 11).Value = 9999.5  # K99
$ws.Cells.Item(99, 13).Value = -8501.5  # M99
$ws.Cells.Item(107, 8).Value = 7169.8335  # H107
$ws.Cells.Item(107, 9).Value = 7124.8237  # I107
$ws.Cells.Item(107, 10).Value = 7279.143  # J107
$ws.Cells.Item(107, 11).Value = 7124.8237  # K107
$ws.Cells.Item(107, 12).Value = 7279.143  # L107
$ws.Cells.Item(107, 13).Value = -5204.8237  # M107
$ws.Cells.Item(107, 14).Value = -11119.143  # N107
$ws.Cells.Item(133, 8).Value = 86233  # H133
$ws.Cells.Item(133, 10).Value = 85495  # J133
$ws.Cells.Item(133, 12).Value = 85495  # L133
$ws.Cells.Item(133, 14).Value = -95615  # N133

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 1930.1904  # H22
$ws.Cells.Item(22, 9).Value = 1825.2727  # I22
$ws.Cells.Item(22, 10).Value = 2045.6  # J22
$ws.Cells.Item(22, 11).Value = 1825.2727  # K22
$ws.Cells.Item(22, 12).Value = 2045.6  # L22
$ws.Cells.Item(22, 13).Value = -1475.2727  # M22
$ws.Cells.Item(22, 14).Value = -2745.6  # N22
$ws.Cells.Item(31, 8).Value = 3243.158  # H31
$ws.Cells.Item(31, 9).Value = 1313.375  # I31
$ws.Cells.Item(31, 10).Value = 4646.636  # J31
$ws.Cells.Item(31, 11).Value = 1313.375  # K31
$ws.Cells.Item(31, 12).Value = 4646.636  # L31
$ws.Cells.Item(31, 13).Value = -1018.375  # M31
$ws.Cells.Item(31, 14).Value = -5236.636  # N31
$ws.Cells.Item(34, 8).Value = 3243.158  # H34
$ws.Cells.Item(34, 9).Value = 1313.375  # I34
$ws.Cells.Item(34, 10).Value = 4646.636  # J34
$ws.Cells.Item(34, 11).Value = 1313.375  # K34
$ws.Cells.Item(34, 12).Value = 4646.636  # L34
$ws.Cells.Item(34, 13).Value = -1111.375  # M34
$ws.Cells.Item(34, 14).Value = -5050.636  # N34
$ws.Cells.Item(58, 8).Value = 5636.6875  # H58
$ws.Cells.Item(58, 9).Value = 5309.4736  # I58
$ws.Cells.Item(58, 10).Value = 6114.923  # J58
$ws.Cells.Item(58, 11).Value = 5309.4736  # K58
$ws.Cells.Item(58, 12).Value = 6114.923  # L58
$ws.Cells.Item(58, 13).Value = -5106.4736  # M58
$ws.Cells.Item(58, 14).Value = -6520.923  # N58
$ws.Cells.Item(122, 8).Value = 4279  # H122
$ws.Cells.Item(122, 9).Value = 3920.1667  # I122
$ws.Cells.Item(122, 10).Value = 4996.6665  # J122
$ws.Cells.Item(122, 11).Value = 11760.5001  # K122
$ws.Cells.Item(122, 12).Value = 14989.9995  # L122
$ws.Cells.Item(122, 13).Value = -9310.500100000001  # M122
$ws.Cells.Item(122, 14).Value = -19889.9995  # N122
$ws.Cells.Item(136, 8).Value = 5636.6875  # H136
$ws.Cells.Item(136, 9).Value = 5309.4736  # I136
$ws.Cells.Item(136, 10).Value = 6114.923  # J136
$ws.Cells.Item(136, 11).Value = 15928.4208  # K136
$ws.Cells.Item(136, 12).Value = 18344.769  # L136
$ws.Cells.Item(136, 13).Value = -13378.4208  # M136
$ws.Cells.Item(136, 14).Value = -23444.769  # N136
$ws.Cells.Item(138, 8).Value = 117925.664  # H138
$ws.Cells.Item(138, 10).Value = 117925.664  # J138
$ws.Cells.Item(138, 12).Value = 117925.664  # L138
$ws.Cells.Item(138, 14).Value = -128205.664  # N138

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 116.27273  # H2
$ws.Cells.Item(2, 9).Value = 63.166668  # I2
$ws.Cells.Item(2, 11).Value = 379.000008  # K2
$ws.Cells.Item(2, 13).Value = -266.000008  # M2
$ws.Cells.Item(3, 8).Value = 2571.111  # H3
$ws.Cells.Item(3, 9).Value = 1642.5  # I3
$ws.Cells.Item(3, 11).Value = 4927.5  # K3
$ws.Cells.Item(3, 13).Value = -4815.5  # M3
$ws.Cells.Item(12, 8).Value = 136  # H12
$ws.Cells.Item(12, 10).Value = 161.3  # J12
$ws.Cells.Item(12, 12).Value = 483.9  # L12
$ws.Cells.Item(12, 14).Value = -829.9000000000001  # N12
$ws.Cells.Item(13, 8).Value = 5237.1113  # H13
$ws.Cells.Item(13, 9).Value = 2050  # I13
$ws.Cells.Item(13, 11).Value = 6150  # K13
$ws.Cells.Item(13, 13).Value = -5982  # M13
$ws.Cells.Item(23, 8).Value = 200.5  # H23
$ws.Cells.Item(23, 9).Value = 166.42857  # I23
$ws.Cells.Item(23, 11).Value = 499.28571  # K23
$ws.Cells.Item(23, 13).Value = -264.28571  # M23
$ws.Cells.Item(38, 8).Value = 379.6  # H38
$ws.Cells.Item(38, 9).Value = 99  # I38
$ws.Cells.Item(38, 10).Value = 566.6667  # J38
$ws.Cells.Item(38, 11).Value = 297  # K38
$ws.Cells.Item(38, 12).Value = 1700.0001  # L38
$ws.Cells.Item(38, 13).Value = 50  # M38
$ws.Cells.Item(38, 14).Value = -2394.0001  # N38
$ws.Cells.Item(116, 8).Value = 8505638  # H116
$ws.Cells.Item(116, 10).Value = 24299  # J116
$ws.Cells.Item(116, 12).Value = 72897  # L116
$ws.Cells.Item(116, 14).Value = -79781  # N116
$ws.Cells.Item(129, 8).Value = 1526.6471  # H129
$ws.Cells.Item(129, 9).Value = 643.1  # I129
$ws.Cells.Item(129, 10).Value = 2788.8572  # J129
$ws.Cells.Item(129, 11).Value = 1929.3  # K129
$ws.Cells.Item(129, 12).Value = 8366.571599999999  # L129
$ws.Cells.Item(129, 13).Value = 3070.7  # M129
$ws.Cells.Item(129, 14).Value = -18366.5716  # N129
$ws.Cells.Item(137, 8).Value = 12868.167  # H137
$ws.Cells.Item(137, 10).Value = 13544.929  # J137
$ws.Cells.Item(137, 12).Value = 40634.787  # L137
$ws.Cells.Item(137, 14).Value = -50834.787  # N137

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(99, 8).Value = 18157.75  # H99
$ws.Cells.Item(99, 9).Value = 16988.555  # I99
$ws.Cells.Item(99, 11).Value = 16988.555  # K99
$ws.Cells.Item(99, 13).Value = -14742.555  # M99
$ws.Cells.Item(132, 8).Value = 6916.5  # H132
$ws.Cells.Item(132, 9).Value = 4000  # I132
$ws.Cells.Item(132, 10).Value = 7499.8  # J132
$ws.Cells.Item(132, 11).Value = 12000  # K132
$ws.Cells.Item(132, 12).Value = 22499.4  # L132
$ws.Cells.Item(132, 13).Value = -9470  # M132
$ws.Cells.Item(132, 14).Value = -27559.4  # N132

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(97, 8).Value = 144804.11  # H97
$ws.Cells.Item(97, 10).Value = 144804.11  # J97
$ws.Cells.Item(97, 12).Value = 144804.11  # L97
$ws.Cells.Item(97, 14).Value = -146786.11  # N97
$ws.Cells.Item(122, 8).Value = 0  # H122
$ws.Cells.Item(122, 9).Value = 0  # I122
$ws.Cells.Item(122, 10).Value = 0  # J122
$ws.Cells.Item(122, 11).Value = 0  # K122
$ws.Cells.Item(122, 12).Value = 0  # L122
$ws.Cells.Item(122, 13).ClearContents() | Out-Null  # M122
$ws.Cells.Item(122, 14).ClearContents() | Out-Null  # N122

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 8907.093999999999  # H122
$ws.Cells.Item(122, 9).Value = 6807.6665  # I122
$ws.Cells.Item(122, 10).Value = 12915.091  # J122
$ws.Cells.Item(122, 11).Value = 20422.9995  # K122
$ws.Cells.Item(122, 12).Value = 38745.273  # L122
$ws.Cells.Item(122, 13).Value = -17972.9995  # M122
$ws.Cells.Item(122, 14).Value = -43645.273  # N122
$ws.Cells.Item(132, 8).Value = 4105.871  # H132
$ws.Cells.Item(132, 9).Value = 3902.9644  # I132
$ws.Cells.Item(132, 11).Value = 11708.8932  # K132
$ws.Cells.Item(132, 13).Value = -9178.893199999999  # M132
$ws.Cells.Item(136, 8).Value = 999  # H136
$ws.Cells.Item(136, 9).Value = 999  # I136
$ws.Cells.Item(136, 11).Value = 2997  # K136
$ws.Cells.Item(136, 13).Value = -447  # M136

Write-Output "Updated cells: 220 set, 2 cleared"
